# Generate Report for Handoff
# Updates the localization-status report: the "zh-cn"/"de-de" status moves
# from "In Translation" to "Ready for handoff", and the corresponding
# handoff timestamps are refreshed. Column widths are re-autofit to
# accommodate the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-03 03:02:47"

# --- zh-cn detail sheet ------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-03 03:02:43"

# --- de-de detail sheet ------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-03 03:02:47"

# --- Resize columns to fit the new, longer status text -----------------
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).AutoFit()
